$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.207114100456238
$ws.Range("B1").Value = 2.426470994949341
$ws.Range("C1").Value = 4.614832878112793
$ws.Range("D1").Value = 2.592143535614014
$ws.Range("E1").Value = 1.094354748725891
